$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "FAPs"
$ws.Range("B2").Value = "Rspo2"
$ws.Range("C2").Value = "Lgr4"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 1.876175666666667
$ws.Range("H2").Value = 5.628527
$ws.Range("I2").Value = 0.9891011365778073
$ws.Range("J2").Value = 0.9927075980877177
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.5779736666666667
$ws.Range("N2").Value = 1.733921
$ws.Range("O2").Value = 0.04835019606981441
$ws.Range("P2").Value = 0.05356228267519154
$ws.Range("Q2").Value = 1.084380129374111
$ws.Range("R2").Value = 9.759421164367
$ws.Range("S2").Value = 0.04782323388641326
$ws.Range("T2").Value = 0.05317168498258477

# Row 3
$ws.Range("A3").Value = "FAPs"
$ws.Range("B3").Value = "Rspo2"
$ws.Range("C3").Value = "Lgr4"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 1.876175666666667
$ws.Range("H3").Value = 5.628527
$ws.Range("I3").Value = 0.9891011365778073
$ws.Range("J3").Value = 0.9927075980877177
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 7.791016
$ws.Range("N3").Value = 23.373048
$ws.Range("O3").Value = 0.6517548686181108
$ws.Range("P3").Value = 0.7220131735856595
$ws.Range("Q3").Value = 14.61731463781067
$ws.Range("R3").Value = 131.555831740296
$ws.Range("S3").Value = 0.6446514813202929
$ws.Range("T3").Value = 0.7167479633379105

# Row 4
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Rspo2"
$ws.Range("C4").Value = "Lgr4"
$ws.Range("D4").Value = "Inflammatory-Mac"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 1.876175666666667
$ws.Range("H4").Value = 5.628527
$ws.Range("I4").Value = 0.9891011365778073
$ws.Range("J4").Value = 0.9927075980877177
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.05015166666666667
$ws.Range("N4").Value = 0.150455
$ws.Range("O4").Value = 0.004195421100317676
$ws.Range("P4").Value = 0.004647681895481942
$ws.Range("Q4").Value = 0.09409333664277779
$ws.Range("R4").Value = 0.8468400297850001
$ws.Range("S4").Value = 0.004149695778746729
$ws.Range("T4").Value = 0.00461378913113965

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Rspo2"
$ws.Range("C5").Value = "Lgr4"
$ws.Range("D5").Value = "MuSCs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 1.876175666666667
$ws.Range("H5").Value = 5.628527
$ws.Range("I5").Value = 0.9891011365778073
$ws.Range("J5").Value = 0.9927075980877177
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 3.489664
$ws.Range("N5").Value = 6.979328
$ws.Range("O5").Value = 0.2919266886169084
$ws.Range("P5").Value = 0.215597330685123
$ws.Range("Q5").Value = 6.547222681642666
$ws.Range("R5").Value = 39.283336089856
$ws.Range("S5").Value = 0.2887450195083797
$ws.Range("T5").Value = 0.2140251082985519

# Row 6
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Rspo2"
$ws.Range("C6").Value = "Lgr4"
$ws.Range("D6").Value = "Resolving-Mac"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 1.876175666666667
$ws.Range("H6").Value = 5.628527
$ws.Range("I6").Value = 0.9891011365778073
$ws.Range("J6").Value = 0.9927075980877177
$ws.Range("K6").Value = 1
$ws.Range("L6").Value = 0.3333333333333333
$ws.Range("M6").Value = 0.0451
$ws.Range("N6").Value = 0.1353
$ws.Range("O6").Value = 0.003772825594848836
$ws.Range("P6").Value = 0.004179531158543795
$ws.Range("Q6").Value = 0.08461552256666667
$ws.Range("R6").Value = 0.7615397031000001
$ws.Range("S6").Value = 0.003731706083974826
$ws.Range("T6").Value = 0.004149052337530787

# Row 7
$ws.Range("A7").Value = "MuSCs"
$ws.Range("B7").Value = "Rspo2"
$ws.Range("C7").Value = "Lgr4"
$ws.Range("D7").Value = "ECs"
$ws.Range("E7").Value = 1
$ws.Range("F7").Value = 0.5
$ws.Range("G7").Value = 0.0206735
$ws.Range("H7").Value = 0.041347
$ws.Range("I7").Value = 0.01089886342219268
$ws.Range("J7").Value = 0.007292401912282354
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 0.5779736666666667
$ws.Range("N7").Value = 1.733921
$ws.Range("O7").Value = 0.04835019606981441
$ws.Range("P7").Value = 0.05356228267519154
$ws.Range("Q7").Value = 0.01194873859783334
$ws.Range("R7").Value = 0.071692431587
$ws.Range("S7").Value = 0.0005269621834011445
$ws.Range("T7").Value = 0.0003905976926067748

# Row 8
$ws.Range("A8").Value = "MuSCs"
$ws.Range("B8").Value = "Rspo2"
$ws.Range("C8").Value = "Lgr4"
$ws.Range("D8").Value = "FAPs"
$ws.Range("E8").Value = 1
$ws.Range("F8").Value = 0.5
$ws.Range("G8").Value = 0.0206735
$ws.Range("H8").Value = 0.041347
$ws.Range("I8").Value = 0.01089886342219268
$ws.Range("J8").Value = 0.007292401912282354
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 7.791016
$ws.Range("N8").Value = 23.373048
$ws.Range("O8").Value = 0.6517548686181108
$ws.Range("P8").Value = 0.7220131735856595
$ws.Range("Q8").Value = 0.161067569276
$ws.Range("R8").Value = 0.966405415656
$ws.Range("S8").Value = 0.007103387297817924
$ws.Range("T8").Value = 0.005265210247749115

# Row 9
$ws.Range("A9").Value = "MuSCs"
$ws.Range("B9").Value = "Rspo2"
$ws.Range("C9").Value = "Lgr4"
$ws.Range("D9").Value = "Inflammatory-Mac"
$ws.Range("E9").Value = 1
$ws.Range("F9").Value = 0.5
$ws.Range("G9").Value = 0.0206735
$ws.Range("H9").Value = 0.041347
$ws.Range("I9").Value = 0.01089886342219268
$ws.Range("J9").Value = 0.007292401912282354
$ws.Range("K9").Value = 1
$ws.Range("L9").Value = 0.3333333333333333
$ws.Range("M9").Value = 0.05015166666666667
$ws.Range("N9").Value = 0.150455
$ws.Range("O9").Value = 0.004195421100317676
$ws.Range("P9").Value = 0.004647681895481942
$ws.Range("Q9").Value = 0.001036810480833333
$ws.Range("R9").Value = 0.006220862885
$ws.Range("S9").Value = 0.00004572532157094769
$ws.Range("T9").Value = 0.00003389276434229259

# Row 10
$ws.Range("A10").Value = "MuSCs"
$ws.Range("B10").Value = "Rspo2"
$ws.Range("C10").Value = "Lgr4"
$ws.Range("D10").Value = "MuSCs"
$ws.Range("E10").Value = 1
$ws.Range("F10").Value = 0.5
$ws.Range("G10").Value = 0.0206735
$ws.Range("H10").Value = 0.041347
$ws.Range("I10").Value = 0.01089886342219268
$ws.Range("J10").Value = 0.007292401912282354
$ws.Range("K10").Value = 2
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 3.489664
$ws.Range("N10").Value = 6.979328
$ws.Range("O10").Value = 0.2919266886169084
$ws.Range("P10").Value = 0.215597330685123
$ws.Range("Q10").Value = 0.072143568704
$ws.Range("R10").Value = 0.288574274816
$ws.Range("S10").Value = 0.003181669108528655
$ws.Range("T10").Value = 0.001572222386571162

# Row 11
$ws.Range("A11").Value = "MuSCs"
$ws.Range("B11").Value = "Rspo2"
$ws.Range("C11").Value = "Lgr4"
$ws.Range("D11").Value = "Resolving-Mac"
$ws.Range("E11").Value = 1
$ws.Range("F11").Value = 0.5
$ws.Range("G11").Value = 0.0206735
$ws.Range("H11").Value = 0.041347
$ws.Range("I11").Value = 0.01089886342219268
$ws.Range("J11").Value = 0.007292401912282354
$ws.Range("K11").Value = 1
$ws.Range("L11").Value = 0.3333333333333333
$ws.Range("M11").Value = 0.0451
$ws.Range("N11").Value = 0.1353
$ws.Range("O11").Value = 0.003772825594848836
$ws.Range("P11").Value = 0.004179531158543795
$ws.Range("Q11").Value = 0.0009323748500000001
$ws.Range("R11").Value = 0.0055942491
$ws.Range("S11").Value = 0.00004111951087401032
$ws.Range("T11").Value = 0.00003047882101300846
